$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy existing formatting from the last existing column (AX) into the
# --- four new columns (AY:BB) before writing values, so the new cells pick
# --- up the same cell styles already used by the table (header style for
# --- row 1, data style for rows 2-6) instead of the default style.
$ws.Range("AX1").Copy()
$ws.Range("AY1:BB1").PasteSpecial(-4122)

$ws.Range("AX2:AX6").Copy()
$ws.Range("AY2:BB6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- New quarter headers (2022-Q1 .. 2022-Q4) ---
$ws.Range("AY1").Value = "2022-Q1"
$ws.Range("AZ1").Value = "2022-Q2"
$ws.Range("BA1").Value = "2022-Q3"
$ws.Range("BB1").Value = "2022-Q4"

# --- China (row 2) ---
$ws.Range("AY2").Value = 0.18
$ws.Range("AZ2").Value = 0.19
$ws.Range("BA2").Value = 0.2
$ws.Range("BB2").Value = 0.24

# --- Japan (row 3) ---
$ws.Range("AY3").Value = 0.18
$ws.Range("AZ3").Value = 0.2
$ws.Range("BA3").Value = 0.2
$ws.Range("BB3").Value = 0.25

# --- EU (row 4) ---
$ws.Range("AY4").Value = 0.18
$ws.Range("AZ4").Value = 0.2
$ws.Range("BA4").Value = 0.2
$ws.Range("BB4").Value = 0.25

# --- USA (row 5) ---
$ws.Range("AY5").Value = 0.27
$ws.Range("AZ5").Value = 0.28000000000000003
$ws.Range("BA5").Value = 0.28999999999999998
$ws.Range("BB5").Value = 0.3

# --- RoW (row 6) ---
$ws.Range("AY6").Value = 0.18
$ws.Range("AZ6").Value = 0.2
$ws.Range("BA6").Value = 0.2
$ws.Range("BB6").Value = 0.24

# --- Restore the view: scrolled over to the new columns, cursor left on BA19 ---
[void]$ws.Range("BA19").Select()
